$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The duplicated line that appears back-to-back (twice) just before the
# "Little Asuka" entries in row 19, and just before "Shiji: Aim the target..."
# in row 115. One copy of each pair needs to be removed (the final, trailing
# occurrence of this line at the very end of each list must be kept).
$dupPair = "'Shinji&Asuka: This is forced by Misato-san, who insists that Japanese should        begin with form.', 'Shinji&Asuka: This is forced by Misato-san, who insists that Japanese should        begin with form.', "

# --- Row 19 (Asuka) ---
$b19 = $ws.Cells.Item(19, 2)
$orig19 = $b19.Text
$new19 = $orig19.Replace($dupPair, "")
if ($new19.Length -eq $orig19.Length) {
    throw "Row 19: duplicate pair not found, no change was made"
}
$b19.Value = $new19
$ws.Cells.Item(19, 3).Value = 173
$ws.Cells.Item(19, 4).Value = 3202

# --- Row 115 (Shinji) ---
$b115 = $ws.Cells.Item(115, 2)
$orig115 = $b115.Text
$new115 = $orig115.Replace($dupPair, "")
if ($new115.Length -eq $orig115.Length) {
    throw "Row 115: duplicate pair not found, no change was made"
}
$b115.Value = $new115
$ws.Cells.Item(115, 3).Value = 235
$ws.Cells.Item(115, 4).Value = 3790
